# This workbook is a weekly/daily price log for "Jengibre" (ginger) at
# "Vega Modelo de Temuco", sorted with the most recent record first
# (row 199) down to the oldest (row 262). The commit adds one new daily
# record at the top of that block (2023-03-03), which pushes every
# existing record in the block down by one row (199->200, ..., 262->263)
# and grows the used range from A1:R262 to A1:R263.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right before the current row 199, shifting
# rows 199:262 down to 200:263 (and the style of D199 carries down with
# them automatically since Excel shifts whole rows).
$ws.Rows.Item(199).Insert()

# Populate the freshly inserted row 199 with the new data point.
$ws.Cells.Item(199, 1).Value = 10
$ws.Cells.Item(199, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(199, 3).Value = "La Araucanía"
$ws.Cells.Item(199, 4).Value = 44988
$ws.Cells.Item(199, 5).Value = 9
$ws.Cells.Item(199, 6).Value = 100114007
$ws.Cells.Item(199, 7).Value = "Jengibre"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 30
$ws.Cells.Item(199, 11).Value = 28000
$ws.Cells.Item(199, 12).Value = 28000
$ws.Cells.Item(199, 13).Value = 28000
$ws.Cells.Item(199, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(199, 15).Value = "Perú"
$ws.Cells.Item(199, 16).Value = 2154
$ws.Cells.Item(199, 17).Value = 13
$ws.Cells.Item(199, 18).Value = "Hortaliza"
